# MOSIP-17570 added supervisor rejected email and sms templates
# Adds RPR_SUP_REJECT_EMAIL / RPR_SUP_REJECT_SMS / RPR_SUP_REJECT_EMAIL_SUBJECT
# template rows (for each supported language) to the template_type sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Languages used throughout the sheet, in the order the new rows were appended.
$langEmailSms = @("eng", "fra", "ara", "hin", "kan", "tam")
$langSubject  = @("eng", "fra", "ara", "hin", "kan", "tam")

$startRow = 1726

# Use the existing last data row (1725) as a style/format template: copy it
# down for every new row so number formats / styles stay identical, then
# overwrite the cell values.
$templateRow = 1725

$row = $startRow

foreach ($lang in $langEmailSms) {
    # RPR_SUP_REJECT_EMAIL row
    $ws.Range("A" + $templateRow + ":D" + $templateRow).Copy() | Out-Null
    $ws.Range("A" + $row + ":D" + $row).PasteSpecial() | Out-Null
    $ws.Range("A" + $row).Value = $lang
    $ws.Range("B" + $row).Value = "RPR_SUP_REJECT_EMAIL"
    $ws.Range("C" + $row).Value = "Template for Supervisor Reject Email"
    $row = $row + 1

    # RPR_SUP_REJECT_SMS row
    $ws.Range("A" + $templateRow + ":D" + $templateRow).Copy() | Out-Null
    $ws.Range("A" + $row + ":D" + $row).PasteSpecial() | Out-Null
    $ws.Range("A" + $row).Value = $lang
    $ws.Range("B" + $row).Value = "RPR_SUP_REJECT_SMS"
    $ws.Range("C" + $row).Value = "Template for Supervisor Reject SMS"
    $row = $row + 1
}

foreach ($lang in $langSubject) {
    # RPR_SUP_REJECT_EMAIL_SUBJECT row
    $ws.Range("A" + $templateRow + ":D" + $templateRow).Copy() | Out-Null
    $ws.Range("A" + $row + ":D" + $row).PasteSpecial() | Out-Null
    $ws.Range("A" + $row).Value = $lang
    $ws.Range("B" + $row).Value = "RPR_SUP_REJECT_EMAIL_SUBJECT"
    $ws.Range("C" + $row).Value = "Template for Supervisor Reject Email Subject"
    $row = $row + 1
}

# Update the view to mirror where the user ended up after adding the rows.
$ws.Range("F1730").Select() | Out-Null
